$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 117 (pushes the existing rows 117-139 down to 118-140,
# and the sheet's used range grows from A1:R139 to A1:R140).
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row with the new weekly price-report entry.
$ws.Cells.Item(117, 1).Value = 4
$ws.Cells.Item(117, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(117, 3).Value = 'Los Lagos'
$ws.Cells.Item(117, 4).Value = 44522
$ws.Cells.Item(117, 5).Value = 10
$ws.Cells.Item(117, 6).Value = 100112028
$ws.Cells.Item(117, 7).Value = 'Sandia'
$ws.Cells.Item(117, 8).Value = 'Sin especificar'
$ws.Cells.Item(117, 9).Value = 'Primera'
$ws.Cells.Item(117, 10).Value = 300
$ws.Cells.Item(117, 11).Value = 1200
$ws.Cells.Item(117, 12).Value = 1200
$ws.Cells.Item(117, 13).Value = 1200
$ws.Cells.Item(117, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(117, 15).Value = 'Perú'
$ws.Cells.Item(117, 16).Value = 1200
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = 'Hortaliza'
